$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.404.82"
$ws.Range("E2").Value = "  -4.15%  "
$ws.Range("D3").Value = "3.094.07"
$ws.Range("E3").Value = "  -4.47%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'548.80"
$ws.Range("E5").Value = "  -4.54%  "
$ws.Range("D6").Value = "'137.43"
$ws.Range("E6").Value = "  -9.85%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "3.086.51"
$ws.Range("E8").Value = "  -4.36%  "
$ws.Range("D9").Value = "'0.498"
$ws.Range("E9").Value = "  -2.73%  "
$ws.Range("D10").Value = "'0.157"
$ws.Range("E10").Value = "  -4.09%  "
$ws.Range("D11").Value = "'6.26"
$ws.Range("E11").Value = "  -11.11%  "
$ws.Range("D12").Value = "'0.471"
$ws.Range("E12").Value = "  -3.06%  "
$ws.Range("D13").Value = "'35.54"
$ws.Range("E13").Value = "  -5.23%  "
$ws.Range("E14").Value = "  -6.61%  "
$ws.Range("D15").Value = "3.594.47"
$ws.Range("E15").Value = "  -4.46%  "
$ws.Range("D16").Value = "63.399.81"
$ws.Range("E16").Value = "  -4.22%  "
$ws.Range("D17").Value = "'0.112"
$ws.Range("E17").Value = "  -2.74%  "
$ws.Range("D18").Value = "3.087.75"
$ws.Range("E18").Value = "  -4.94%  "
$ws.Range("D19").Value = "'6.76"
$ws.Range("E19").Value = "  -4.40%  "
$ws.Range("D20").Value = "'490.01"
$ws.Range("E20").Value = "  -11.76%  "
$ws.Range("E21").Value = "  -4.85%  "
$ws.Range("D22").Value = "'0.719"
$ws.Range("E22").Value = "  -2.76%  "
$ws.Range("D23").Value = "'7.27"
$ws.Range("E23").Value = "  -5.88%  "
$ws.Range("D24").Value = "'79.05"
$ws.Range("E24").Value = "  -3.08%  "
$ws.Range("D25").Value = "'12.39"
$ws.Range("E25").Value = "  -8.26%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").Value = "'8.51"
$ws.Range("E27").Value = "  -7.79%  "
$ws.Range("D28").Value = "'2.76"
$ws.Range("E28").Value = "  -5.78%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "'1.98"
$ws.Range("E29").Value = "  -10.64%  "
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").Value = "'26.68"
$ws.Range("E31").Value = "  -3.43%  "
$ws.Range("E32").Value = "  -4.10%  "
$ws.Range("D33").Value = "'2.51"
$ws.Range("E33").Value = "  -7.61%  "
$ws.Range("D34").Value = "'58.60"
$ws.Range("E34").Value = "  +6.06%  "
$ws.Range("D35").Value = "'512.07"
$ws.Range("E35").Value = "  -8.48%  "
$ws.Range("D36").Value = "'6.06"
$ws.Range("E36").Value = "  -4.38%  "
$ws.Range("D37").Value = "'5.14"
$ws.Range("E37").Value = "  -9.55%  "
$ws.Range("D38").Value = "'0.0401"
$ws.Range("E38").Value = "  -11.35%  "
$ws.Range("D39").Value = "3.154.81"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").Value = "'0.0804"
$ws.Range("E40").Value = "  -6.23%  "
$ws.Range("D41").Value = "'0.120"
$ws.Range("E41").Value = "  -7.73%  "
$ws.Range("D42").Value = "'8.18"
$ws.Range("E42").Value = "  -4.42%  "
$ws.Range("E43").Value = "  -12.18%  "
$ws.Range("E44").Value = "  -4.95%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'2.08"
$ws.Range("E46").Value = "  -8.65%  "
$ws.Range("D47").Value = "'25.31"
$ws.Range("E47").Value = "  -3.84%  "
$ws.Range("D48").Value = "'121.02"
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("E49").Value = "  -2.71%  "
$ws.Range("D50").Value = "0.0₃0506"
$ws.Range("E50").Value = "  -7.84%  "
$ws.Range("D51").Value = "'2.32"
$ws.Range("E51").Value = "  +33.12%  "
